# ------------------------------------------------------------------
# Update "北京-漫展信息" workbook to the commit's published data.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "展览" (exhibitions) - 想去人数 (F column) refreshes
# -----------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    3  = 3199
    4  = 25
    5  = 1343
    7  = 385
    8  = 187
    10 = 8326
    11 = 447
    13 = 71
    14 = 270
    15 = 302
    18 = 336
    19 = 10480
    23 = 24
    26 = 387
    27 = 174
    28 = 153
    30 = 43
    32 = 2069
    33 = 34
    35 = 869
    37 = 271
    40 = 1226
    41 = 161
    44 = 321
    45 = 264
    50 = 64
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# -----------------------------------------------------------------
# Sheet "演出" (performances) - refresh two counts, then insert a
# brand-new show (2024-06-01) above the existing row 19, pushing the
# following four rows down by one.
# -----------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Cells.Item(14, 6).Value = 32
$wsShow.Cells.Item(18, 6).Value = 173

# Insert a new row 19; carry the formatting down from the row that is
# about to become row 20 (the previous row 19) so the new row keeps the
# same cell style (bordered/bold index column etc.).
$wsShow.Rows.Item(19).Insert()
$wsShow.Cells.Item(20, 1).Copy()
$wsShow.Cells.Item(19, 1).PasteSpecial(-4122)

$wsShow.Cells.Item(19, 1).Value = 18
$wsShow.Cells.Item(19, 2).Value = "'2024-06-01"
$wsShow.Cells.Item(19, 3).Value = "北京·六一特献｜【直到世界尽头】灌篮高手等神级中日动漫演唱会，全体起立！"
$wsShow.Cells.Item(19, 4).Value = "学清路38号金码大厦B座 北京想象空间"
$wsShow.Cells.Item(19, 5).Value = "2024.06.01 20:00-06.01 22:00"
$wsShow.Cells.Item(19, 6).Value = 0
$wsShow.Cells.Item(19, 7).Value = 128
$wsShow.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84418"
$wsShow.Cells.Item(19, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/T1N8OnSu1713252809650.jpeg"

# -----------------------------------------------------------------
# Sheet "本地生活" (local life) - single 想去人数 refresh
# -----------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Cells.Item(3, 6).Value = 10

# -----------------------------------------------------------------
# Sheet "全部类型" (all types, merged view) - 想去人数 refreshes
# -----------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 3199
    4  = 25
    5  = 1343
    6  = 385
    9  = 187
    11 = 8326
    12 = 447
    15 = 71
    16 = 270
    17 = 302
    19 = 336
    20 = 10481
    26 = 387
    27 = 174
    28 = 32
    29 = 153
    31 = 2069
    32 = 34
    34 = 869
    35 = 173
    36 = 271
    39 = 1226
    40 = 161
    43 = 321
    45 = 264
    50 = 64
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
